$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$cD.NumberFormat = "@"
$cD.Value2 = "63.298.92"
$cD.Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "

$cD = $ws.Range("D3")
$cD.NumberFormat = "@"
$cD.Value2 = "2.648.32"
$cD.Style = "Normal"
$ws.Range("E3").Value = "  +2.82%  "

$cD = $ws.Range("D4")
$cD.NumberFormat = "@"
$cD.Value2 = "0.999"
$cD.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

$cD = $ws.Range("D5")
$cD.NumberFormat = "@"
$cD.Value2 = "607.28"
$cD.Style = "Normal"
$ws.Range("E5").Value = "  +3.49%  "

$cD = $ws.Range("D6")
$cD.NumberFormat = "@"
$cD.Value2 = "144.43"
$cD.Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "

$cD = $ws.Range("D7")
$cD.NumberFormat = "@"
$cD.Value2 = "0.999"
$cD.Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  -0.27%  "

$cD = $ws.Range("D9")
$cD.NumberFormat = "@"
$cD.Value2 = "2.647.39"
$cD.Style = "Normal"
$ws.Range("E9").Value = "  +2.87%  "

$ws.Range("E10").Value = "  +1.56%  "

$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("E13").Value = "  +4.24%  "

$cD = $ws.Range("D14")
$cD.NumberFormat = "@"
$cD.Value2 = "27.38"
$cD.Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "

$cD = $ws.Range("D15")
$cD.NumberFormat = "@"
$cD.Value2 = "3.120.17"
$cD.Style = "Normal"
$ws.Range("E15").Value = "  +2.77%  "

$cD = $ws.Range("D16")
$cD.NumberFormat = "@"
$cD.Value2 = "63.097.51"
$cD.Style = "Normal"
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("E17").Value = "  -0.45%  "

$cD = $ws.Range("D18")
$cD.NumberFormat = "@"
$cD.Value2 = "2.628.98"
$cD.Style = "Normal"
$ws.Range("E18").Value = "  +1.91%  "

$cD = $ws.Range("D19")
$cD.NumberFormat = "@"
$cD.Value2 = "11.40"
$cD.Style = "Normal"
$ws.Range("E19").Value = "  +2.60%  "

$cD = $ws.Range("D20")
$cD.NumberFormat = "@"
$cD.Value2 = "343.41"
$cD.Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("E21").Value = "  +3.03%  "

$ws.Range("E22").Value = "  +3.50%  "

$ws.Range("E23").Value = "  -0.11%  "

$cD = $ws.Range("D24")
$cD.NumberFormat = "@"
$cD.Value2 = "67.14"
$cD.Style = "Normal"
$ws.Range("E24").Value = "  -1.10%  "

$cD = $ws.Range("D25")
$cD.NumberFormat = "@"
$cD.Value2 = "1.65"
$cD.Style = "Normal"
$ws.Range("E25").Value = "  +1.27%  "

$ws.Range("E26").Value = "  -2.01%  "

$cD = $ws.Range("D27")
$cD.NumberFormat = "@"
$cD.Value2 = "8.67"
$cD.Style = "Normal"
$ws.Range("E27").Value = "  +5.49%  "

$cD = $ws.Range("D28")
$cD.NumberFormat = "@"
$cD.Value2 = "0.165"
$cD.Style = "Normal"
$ws.Range("E28").Value = "  +0.26%  "

$cD = $ws.Range("D29")
$cD.NumberFormat = "@"
$cD.Value2 = "547.18"
$cD.Style = "Normal"
$ws.Range("E29").Value = "  +15.74%  "

$ws.Range("E30").Value = "  +0.24%  "

$cD = $ws.Range("D31")
$cD.NumberFormat = "@"
$cD.Value2 = "7.92"
$cD.Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("E32").Value = "  +4.81%  "

$ws.Range("E33").Value = "  +6.96%  "

$cD = $ws.Range("D34")
$cD.NumberFormat = "@"
$cD.Value2 = "0.0₃0813"
$cD.Style = "Normal"
$ws.Range("E34").Value = "  +1.76%  "

$cD = $ws.Range("D35")
$cD.NumberFormat = "@"
$cD.Value2 = "172.55"
$cD.Style = "Normal"
$ws.Range("E35").Value = "  -2.11%  "

$ws.Range("E36").Value = "  +12.00%  "

$ws.Range("E37").Value = "  +2.13%  "

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("E39").Value = "  +1.60%  "

$cD = $ws.Range("D40")
$cD.NumberFormat = "@"
$cD.Value2 = "1.84"
$cD.Style = "Normal"
$ws.Range("E40").Value = "  +6.44%  "

$cD = $ws.Range("D41")
$cD.NumberFormat = "@"
$cD.Value2 = "171.72"
$cD.Style = "Normal"
$ws.Range("E41").Value = "  +8.20%  "

$cD = $ws.Range("D42")
$cD.NumberFormat = "@"
$cD.Value2 = "0.999"
$cD.Style = "Normal"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("E43").Value = "  +1.49%  "

$cD = $ws.Range("D44")
$cD.NumberFormat = "@"
$cD.Value2 = "22.39"
$cD.Style = "Normal"
$ws.Range("E44").Value = "  +3.67%  "

$ws.Range("E45").Value = "  +7.52%  "

$cD = $ws.Range("D46")
$cD.NumberFormat = "@"
$cD.Value2 = "0.632"
$cD.Style = "Normal"
$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cD = $ws.Range("D47")
$cD.NumberFormat = "@"
$cD.Value2 = "0.0240"
$cD.Style = "Normal"
$ws.Range("E47").Value = "  +1.56%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cD = $ws.Range("D48")
$cD.NumberFormat = "@"
$cD.Value2 = "0.0962"
$cD.Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "

$cD = $ws.Range("D49")
$cD.NumberFormat = "@"
$cD.Value2 = "18.86"
$cD.Style = "Normal"
$ws.Range("E49").Value = "  +4.15%  "

$ws.Range("E50").Value = "  +3.51%  "

$ws.Range("E51").Value = "  -1.32%  "
